$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.224.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4663"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2712"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06277"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.841.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07413"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.931"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6195"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.158.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007279"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -4.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.885"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.851"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.203"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.865"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1034"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.090"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.809"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04823"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7086"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.695"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.650"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8932"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -6.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.544"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.034"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1194"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.555"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.358"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3641"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.03%  "
